# Update "想去人数" (interested-count) values for two events that appear
# on multiple sheets ("展览" and "全部类型").
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): row3 -> 南宁·熊喵M动漫嘉年华·万圣派对 (156 -> 160)
#                             row4 -> 南宁·万圣漫控嘉年华10          (723 -> 727)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 160
$wsExhibit.Range("F4").Value = 727

# Sheet "全部类型" (all categories): same two events appear one row lower
# row4 -> 南宁·熊喵M动漫嘉年华·万圣派对 (156 -> 160)
# row5 -> 南宁·万圣漫控嘉年华10          (723 -> 727)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 160
$wsAll.Range("F5").Value = 727
